$d = $word.ActiveDocument

# Shared run/paragraph-mark formatting used throughout the document:
# majorHAnsi theme font at 7.5pt (sz/szCs = 15 half-points).
$rpr = '<w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="15"/><w:szCs w:val="15"/></w:rPr>'

# Three new paragraphs to insert at the very top of the document body:
#   1. "Verified that the below remains accurate:"
#   2. <tab>12/17/23 - DHB   (tab + date in one run, " - DHB" in a second run)
#   3. (empty paragraph)
$body =
  '<w:p><w:pPr>' + $rpr + '</w:pPr>' +
    '<w:r>' + $rpr + '<w:t>Verified that the below remains accurate:</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr>' + $rpr + '</w:pPr>' +
    '<w:r>' + $rpr + '<w:tab/><w:t>12/17/23</w:t></w:r>' +
    '<w:r>' + $rpr + '<w:t xml:space="preserve"> - DHB</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr>' + $rpr + '</w:pPr></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $body + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

# Insert the new paragraphs before the first paragraph of the document.
$r = $d.Range(0, 0)
$null = $r.InsertXML($xml)
